# Update "想去人数" (column F) counts on the 展览, 本地生活 and 全部类型
# sheets to reflect the latest gh-pages generation snapshot.

$wb = $excel.ActiveWorkbook

$exhibitionWs = $wb.Worksheets.Item("展览")
$localLifeWs  = $wb.Worksheets.Item("本地生活")
$allTypesWs   = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1) - column F updates
$exhibitionWs.Range("F4").Value  = 8096
$exhibitionWs.Range("F11").Value = 533
$exhibitionWs.Range("F21").Value = 11497
$exhibitionWs.Range("F22").Value = 118
$exhibitionWs.Range("F23").Value = 2227
$exhibitionWs.Range("F25").Value = 3077
$exhibitionWs.Range("F28").Value = 2648
$exhibitionWs.Range("F31").Value = 274
$exhibitionWs.Range("F34").Value = 2351
$exhibitionWs.Range("F39").Value = 5769
$exhibitionWs.Range("F42").Value = 826

# 本地生活 (sheet3) - column F updates
$localLifeWs.Range("F2").Value = 213
$localLifeWs.Range("F3").Value = 351

# 全部类型 (sheet4) - column F updates
$allTypesWs.Range("F4").Value  = 213
$allTypesWs.Range("F5").Value  = 351
$allTypesWs.Range("F8").Value  = 8096
$allTypesWs.Range("F14").Value = 533
$allTypesWs.Range("F24").Value = 11497
$allTypesWs.Range("F25").Value = 118
$allTypesWs.Range("F26").Value = 2227
$allTypesWs.Range("F27").Value = 2227
$allTypesWs.Range("F28").Value = 3077
$allTypesWs.Range("F29").Value = 2648
$allTypesWs.Range("F31").Value = 274
$allTypesWs.Range("F34").Value = 2351
$allTypesWs.Range("F39").Value = 5769
$allTypesWs.Range("F44").Value = 826
